# Auto-generated Excel COM-interop script
# Applies numeric updates to the Famfrit_Profits leve-profit sheets
# per the scheduled price-refresh runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4544.8213
$ws.Range("J38").Value = 5606.5884
$ws.Range("L38").Value = 16819.7652
$ws.Range("N38").Value = -17563.7652

$ws.Range("H62").Value = 1315.5
$ws.Range("I62").Value = 1315.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1315.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -691.5
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 1315.5
$ws.Range("I65").Value = 1315.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 6577.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -3457.5
$ws.Range("N65").ClearContents()

$ws.Range("H74").Value = 5669.0586
$ws.Range("I74").Value = 5166.6665
$ws.Range("J74").Value = 5943.091
$ws.Range("K74").Value = 5166.6665
$ws.Range("L74").Value = 5943.091
$ws.Range("M74").Value = -4230.6665
$ws.Range("N74").Value = -7815.091

$ws.Range("H77").Value = 5669.0586
$ws.Range("I77").Value = 5166.6665
$ws.Range("J77").Value = 5943.091
$ws.Range("K77").Value = 25833.3325
$ws.Range("L77").Value = 29715.455
$ws.Range("M77").Value = -21153.3325
$ws.Range("N77").Value = -39075.455

$ws.Range("H113").Value = 6208.7
$ws.Range("I113").Value = 5579.8
$ws.Range("K113").Value = 5579.8
$ws.Range("M113").Value = -2325.8

$ws.Range("H135").Value = 623.64703
$ws.Range("I135").Value = 568.875
$ws.Range("K135").Value = 5119.875
$ws.Range("M135").Value = -2584.875

$ws.Range("H137").Value = 12435.514
$ws.Range("I137").Value = 12280.866
$ws.Range("K137").Value = 36842.598
$ws.Range("M137").Value = -34292.598

$ws.Range("H138").Value = 24394562
$ws.Range("I138").Value = 875.24
$ws.Range("J138").Value = 62509696
$ws.Range("K138").Value = 2625.72
$ws.Range("L138").Value = 187529088
$ws.Range("M138").Value = 2514.28
$ws.Range("N138").Value = -187539368

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2883.5
$ws.Range("I2").Value = 2584.4
$ws.Range("K2").Value = 2584.4
$ws.Range("M2").Value = -2471.4

$ws.Range("H32").Value = 5122.125
$ws.Range("I32").Value = 5221.243
$ws.Range("K32").Value = 5221.243
$ws.Range("M32").Value = -4934.243

$ws.Range("H74").Value = 297029.75
$ws.Range("I74").Value = 297029.75
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 297029.75
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -296155.75
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 297029.75
$ws.Range("I77").Value = 297029.75
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 1485148.75
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1480780.75
$ws.Range("N77").ClearContents()

$ws.Range("H97").Value = 1611.4736
$ws.Range("I97").Value = 1319.2142
$ws.Range("K97").Value = 1319.2142
$ws.Range("M97").Value = -823.2141999999999

$ws.Range("H116").Value = 2883.5
$ws.Range("I116").Value = 2584.4
$ws.Range("K116").Value = 2584.4
$ws.Range("M116").Value = -290.4000000000001

$ws.Range("H132").Value = 54206.562
$ws.Range("I132").Value = 11331.379
$ws.Range("J132").Value = 468666.66
$ws.Range("K132").Value = 33994.137
$ws.Range("L132").Value = 1405999.98
$ws.Range("M132").Value = -31464.137
$ws.Range("N132").Value = -1411059.98

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2883.5
$ws.Range("I3").Value = 2584.4
$ws.Range("K3").Value = 2584.4
$ws.Range("M3").Value = -2470.4

$ws.Range("H20").Value = 1110.7742
$ws.Range("I20").Value = 1003.36
$ws.Range("J20").Value = 1558.3334
$ws.Range("K20").Value = 1003.36
$ws.Range("L20").Value = 1558.3334
$ws.Range("M20").Value = -756.36
$ws.Range("N20").Value = -2052.3334

$ws.Range("H34").Value = 19963.5
$ws.Range("J34").Value = 19963.5
$ws.Range("L34").Value = 19963.5
$ws.Range("N34").Value = -20191.5

$ws.Range("H94").Value = 1575.8334
$ws.Range("I94").Value = 908.5
$ws.Range("J94").Value = 2576.8333
$ws.Range("K94").Value = 908.5
$ws.Range("L94").Value = 2576.8333
$ws.Range("M94").Value = -457.5
$ws.Range("N94").Value = -3478.8333

$ws.Range("H105").Value = 9424.316999999999
$ws.Range("I105").Value = 20730.455
$ws.Range("K105").Value = 20730.455
$ws.Range("M105").Value = -18983.455

$ws.Range("H113").Value = 3700
$ws.Range("I113").Value = 3700
$ws.Range("K113").Value = 3700
$ws.Range("M113").Value = -1530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 6998.3335
$ws.Range("I22").Value = 6998.3335
$ws.Range("K22").Value = 6998.3335
$ws.Range("M22").Value = -6648.3335

$ws.Range("H86").Value = 9123.75
$ws.Range("I86").Value = 7248
$ws.Range("K86").Value = 7248
$ws.Range("M86").Value = -6125

$ws.Range("H89").Value = 9123.75
$ws.Range("I89").Value = 7248
$ws.Range("K89").Value = 36240
$ws.Range("M89").Value = -30624

$ws.Range("H134").Value = 1555.3529
$ws.Range("I134").Value = 1309.4667
$ws.Range("J134").Value = 3399.5
$ws.Range("K134").Value = 3928.4001
$ws.Range("L134").Value = 10198.5
$ws.Range("M134").Value = -1393.4001
$ws.Range("N134").Value = -15268.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 193.3
$ws.Range("I2").Value = 38.6
$ws.Range("J2").Value = 348
$ws.Range("K2").Value = 38.6
$ws.Range("L2").Value = 348
$ws.Range("M2").Value = 74.40000000000001
$ws.Range("N2").Value = -574

$ws.Range("H17").Value = 450
$ws.Range("I17").Value = 100
$ws.Range("K17").Value = 100
$ws.Range("M17").Value = 68

$ws.Range("H32").Value = 53199.2
$ws.Range("I32").Value = 41998.668
$ws.Range("K32").Value = 41998.668
$ws.Range("M32").Value = -41702.668

$ws.Range("H41").Value = 13493.75
$ws.Range("I41").Value = 7966.6665
$ws.Range("J41").Value = 30075
$ws.Range("K41").Value = 7966.6665
$ws.Range("L41").Value = 30075
$ws.Range("M41").Value = -7611.6665
$ws.Range("N41").Value = -30785

$ws.Range("H42").Value = 52500
$ws.Range("I42").Value = 52500
$ws.Range("K42").Value = 52500
$ws.Range("M42").Value = -52015

$ws.Range("H70").Value = 120309.664
$ws.Range("I70").Value = 132629.81
$ws.Range("J70").Value = 21748.5
$ws.Range("K70").Value = 132629.81
$ws.Range("L70").Value = 21748.5
$ws.Range("M70").Value = -132359.81
$ws.Range("N70").Value = -22288.5

$ws.Range("H73").Value = 120309.664
$ws.Range("I73").Value = 132629.81
$ws.Range("J73").Value = 21748.5
$ws.Range("K73").Value = 132629.81
$ws.Range("L73").Value = 21748.5
$ws.Range("M73").Value = -131693.81
$ws.Range("N73").Value = -23620.5

$ws.Range("H115").Value = 52500
$ws.Range("I115").Value = 52500
$ws.Range("K115").Value = 52500
$ws.Range("M115").Value = -51325

$ws.Range("H123").Value = 56710.715
$ws.Range("J123").Value = 58744
$ws.Range("L123").Value = 58744
$ws.Range("N123").Value = -63644

$ws.Range("H140").Value = 105085.89
$ws.Range("J140").Value = 105085.89
$ws.Range("L140").Value = 105085.89
$ws.Range("N140").Value = -115445.89

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5452
$ws.Range("M20").ClearContents()

$ws.Range("H43").Value = 22138.908
$ws.Range("I43").Value = 20999.5
$ws.Range("K43").Value = 20999.5
$ws.Range("M43").Value = -20806.5

$ws.Range("H82").Value = 2665.3215
$ws.Range("I82").Value = 2583.5
$ws.Range("J82").Value = 2774.4167
$ws.Range("K82").Value = 2583.5
$ws.Range("L82").Value = 2774.4167
$ws.Range("M82").Value = -2222.5
$ws.Range("N82").Value = -3496.4167

$ws.Range("H85").Value = 2665.3215
$ws.Range("I85").Value = 2583.5
$ws.Range("J85").Value = 2774.4167
$ws.Range("K85").Value = 2583.5
$ws.Range("L85").Value = 2774.4167
$ws.Range("M85").Value = -1335.5
$ws.Range("N85").Value = -5270.4167

$ws.Range("H122").Value = 3426.8333
$ws.Range("I122").Value = 2986.3157
$ws.Range("J122").Value = 5100.8
$ws.Range("K122").Value = 8958.947100000001
$ws.Range("L122").Value = 15302.4
$ws.Range("M122").Value = -6508.947100000001
$ws.Range("N122").Value = -20202.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1322.7858
$ws.Range("I107").Value = 1525.1666
$ws.Range("K107").Value = 4575.4998
$ws.Range("M107").Value = -2655.4998

$ws.Range("H122").Value = 55070.15
$ws.Range("I122").Value = 87012.914
$ws.Range("J122").Value = 7156
$ws.Range("K122").Value = 261038.742
$ws.Range("L122").Value = 21468
$ws.Range("M122").Value = -258588.742
$ws.Range("N122").Value = -26368
